# Fixed errors with merge
# Insert a new row for "CatalogEntry.estimatedDuration" above the existing
# "CatalogEntry.billingCode" row (old row 26), shifting billingCode,
# billingSummary, scheduleSummary, limitationSummary and regulatorySummary
# down by one row. Also fix two pre-existing data errors that were
# uncovered by the merge: billingSummary.Is Summary? should be "Y" and
# scheduleSummary's RIM mapping should be "OM1-40".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 26 (pushes old rows 26-30 down to 27-31)
$ws.Rows.Item(26).Insert()

# Match the formatting (style "2") used by all the other data rows
$ws.Range("A27:AK27").Copy()
$ws.Range("A26:AK26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 26 with the CatalogEntry.estimatedDuration entry
$ws.Cells.Item(26, 1).Value = "CatalogEntry.estimatedDuration"   # A Path
$ws.Cells.Item(26, 5).Value = "0"                                  # E Min
$ws.Cells.Item(26, 6).Value = "1"                                  # F Max
$ws.Cells.Item(26, 10).Value = "Duration`n"                        # J Type(s)
$ws.Cells.Item(26, 11).Value = "Estimated duration of the orderable item"  # K Short
$ws.Cells.Item(26, 12).Value = "Estimated duration of the orderable item of this  entry of the catalog."  # L Definition
$ws.Cells.Item(26, 31).Value = "CatalogEntry.estimatedDuration"   # AE Base Path
$ws.Cells.Item(26, 32).Value = "0"                                  # AF Base Min
$ws.Cells.Item(26, 33).Value = "1"                                  # AG Base Max

# Fix the billingSummary row (now row 28): "Is Summary?" should be Y
$ws.Cells.Item(28, 9).Value = "Y"

# Fix the scheduleSummary row (now row 29): RIM mapping should be OM1-40
$ws.Cells.Item(29, 36).Value = "OM1-40"
